$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.639683365821838
$ws.Range("B1").Value = 3.312361240386963
$ws.Range("C1").Value = 4.223691463470459
$ws.Range("D1").Value = 1.31499981880188
$ws.Range("E1").Value = 0.7718048095703125
